$d = $word.ActiveDocument

# 1. "In fact, a recent" -> "A recent"
$d.Content.Find.Execute("In fact, a recent", $true, $false, $false, $false, $false, $true, 1, $false, "A recent", 2) | Out-Null

# 2. "Mckinsey [1], concluded" -> "Mckinsey, concluded"
$d.Content.Find.Execute("Mckinsey [1], concluded", $true, $false, $false, $false, $false, $true, 1, $false, "Mckinsey, concluded", 2) | Out-Null

# 3. "overall economic productivity, R&D spending" -> "overall economic productivity. R&D spending"
$d.Content.Find.Execute("overall economic productivity, R&D spending", $true, $false, $false, $false, $false, $true, 1, $false, "overall economic productivity. R&D spending", 2) | Out-Null

# 4. "These new technologies, will disrupt" -> "These new technologies will disrupt"
$d.Content.Find.Execute("These new technologies, will disrupt", $true, $false, $false, $false, $false, $true, 1, $false, "These new technologies will disrupt", 2) | Out-Null

# 5. "workflows and in the medium term the role" -> "workflows and, in the medium term, the role"
$d.Content.Find.Execute("workflows and in the medium term the role", $true, $false, $false, $false, $false, $true, 1, $false, "workflows and, in the medium term, the role", 2) | Out-Null

# 6. "Currently most of these tools" -> "Currently, most of these tools"
$d.Content.Find.Execute("Currently most of these tools", $true, $false, $false, $false, $false, $true, 1, $false, "Currently, most of these tools", 2) | Out-Null

# 7. "but also to the final stage" -> "but also in the final stage"
$d.Content.Find.Execute("but also to the final stage", $true, $false, $false, $false, $false, $true, 1, $false, "but also in the final stage", 2) | Out-Null

# 8. "when compared to other options (buildings" -> "in comparison to other types of structures (buildings"
$d.Content.Find.Execute("when compared to other options (buildings", $true, $false, $false, $false, $false, $true, 1, $false, "in comparison to other types of structures (buildings", 2) | Out-Null

# 9. "material savings when compared to the case study" -> "material savings in comparison to the case study"
$d.Content.Find.Execute("material savings when compared to the case study", $true, $false, $false, $false, $false, $true, 1, $false, "material savings in comparison to the case study", 2) | Out-Null

Write-Output "text edits done"

# 10. Move the _GoBack bookmark from inside "disruption" (after "disr") to
#     inside "Optimization" (after "Optimi") in the Keywords paragraph, and
#     split that run into "Optimi" / "zation, Steel Tower, Lattice, Electricity pylon".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$kwRange = $d.Content
$kwRange.Find.Execute("Optimization, Steel Tower", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $kwRange.Start + 6   # after "Optimi"
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 11. Delete the whole separator-line paragraph
#     "------------------------------------- -- -------------------------"
$sepRange = $d.Content
$sepRange.Find.Execute("------------------------------------- -- -------------------------", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sepPara = $sepRange.Paragraphs(1)
$sepPara.Range.Delete()

# 12. Remove the footnote text "[1] - Mckinsey&Company, 2016 , ..." but keep
#     the (now empty) paragraph it lived in (no leftover empty run).
$fnRange = $d.Content
$fnRange.Find.Execute("[1] " + [char]8211 + " Mckinsey&Company", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fnPara = $fnRange.Paragraphs(1)
$fnParaRange = $fnPara.Range
$fnParaRange.MoveEnd(1, -1) | Out-Null
$fnParaRange.Delete()

Write-Output "structural edits done"

